# Update "Esperado" (C), "Observado" (D) and "valor p" (E) columns
# with the figures for semana 28 y 29 de 2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r3 - evento 113 (Desnutricion aguda en menores de 5 anos)
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.37

# r5 - evento 155 (Cancer de la mama y cuello uterino)
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 0.14

# r6 - evento 210 (Dengue)
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 123

# r7 - evento 215 (Defectos congenitos)
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0.14

# r11 - evento 300 (Agresiones por animales potencialmente transmisores de rabia)
$ws.Range("C11").Value = 46
$ws.Range("D11").Value = 21
$ws.Range("E11").Value = 0

# r13 - evento 340 (Hepatitis b, c y coinfeccion hepatitis b y delta)
$ws.Range("C13").Value = 2
$ws.Range("E13").Value = 0.14

# r15 - evento 346 (Ira por virus nuevo)
$ws.Range("C15").Value = 58
$ws.Range("D15").Value = 2

# r16 - evento 348 (Infeccion respiratoria aguda grave irag inusitada)
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0.01

# r17 - evento 352 (Infecciones de sitio quirurgico asociadas a procedimiento medico quirurgico)
$ws.Range("C17").Value = 0
$ws.Range("E17").Value = 1

# r18 - evento 355 (Enfermedad transmitida por alimentos o agua (eta))
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 0.14

# r19 - evento 356 (Intento de suicidio)
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 1

# r20 - evento 357 (Iad - infecciones asociadas a dispositivos - individual)
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 0.02

# r21 - evento 365 (Intoxicaciones)
$ws.Range("C21").Value = 7
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = 0.05

# r24 - evento 455 (Leptospirosis)
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 0.18

# r27 - evento 549 (Morbilidad materna extrema)
$ws.Range("C27").Value = 7
$ws.Range("D27").Value = 7
$ws.Range("E27").Value = 0.15

# r28 - evento 560 (Mortalidad perinatal y neonatal tardia)
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0.37

# r31 - evento 740 (Sifilis congenita)
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 1

# r32 - evento 750 (Sifilis gestacional)
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = 0

# r33 - evento 813 (Tuberculosis)
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0

# r34 - evento 831 (Varicela individual)
$ws.Range("C34").Value = 11
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = 0

# r35 - evento 850 (Vih/sida/mortalidad por sida)
$ws.Range("C35").Value = 8
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
